# Configurations done for sanity suite to be run on prod
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Rename TCID column values from TestCase_F1..F21 to Notifications001..021
for ($i = 1; $i -le 21; $i++) {
    $row = $i + 1
    $num = "{0:D3}" -f $i
    $ws.Cells.Item($row, 1).Value = "Notifications$num"
}

# Flip Runmode column (D) from "N" to "Y" for rows 4 through 22
for ($row = 4; $row -le 22; $row++) {
    $ws.Cells.Item($row, 4).Value = "Y"
}

# Widen column A slightly
$ws.Columns.Item(1).ColumnWidth = 16.140625

# Update the active selection / scroll position
$ws.Range("C12").Select()
